$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2, shifting all existing data down.
$ws.Rows("2:2").Insert()

# Populate the new row with the new Eurobarometer wave entry (ZA7846 / EB 96.1).
$ws.Range("A2").Value = "ZA7846"
$ws.Range("B2").Value = "'96.1"
$ws.Range("C2").Value = "September-October 2021"
$ws.Range("D2").Value = "Future of Europe, and Digital rights and principles (COVID-19 Pandemic)"

# Match the author's final selection recorded in the saved file.
$ws.Range("D3").Select() | Out-Null
